$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.25132410141365114
$ws.Range("B1").Value = 0.25110367288657898
$ws.Range("A2").Value = -0.22900203886646153
$ws.Range("B2").Value = 0.22809334693160821
$ws.Range("A3").Value = -0.13414531606837699
$ws.Range("B3").Value = 0.13389476335503225
$ws.Range("A4").Value = -0.12589476345188189
$ws.Range("B4").Value = 0.12507651544148679
$ws.Range("A5").Value = -0.084011329795849576
$ws.Range("B5").Value = 0.082804369681296563
$ws.Range("A6").Value = -0.039187556606128382
$ws.Range("B6").Value = 0.038835326932463943
$ws.Range("A7").Value = -0.028835327067985084
$ws.Range("B7").Value = 0.028758500994626068
$ws.Range("A8").Value = -0.018758501133382399
$ws.Range("B8").Value = 0.018645447286584016
$ws.Range("A9").Value = -0.016645447360402965
$ws.Range("B9").Value = 0.016560305715823009
$ws.Range("A10").Value = -0.014560305792738149
$ws.Range("B10").Value = 0.014557054426099114
$ws.Range("A11").Value = -0.011557054511554199
$ws.Range("B11").Value = 0.011550889058934644
$ws.Range("A12").Value = -0.008050889149281204
$ws.Range("B12").Value = 0.0080101760407047529
$ws.Range("A13").Value = -0.0045101761332810852
$ws.Range("B13").Value = 0.0044985009519882269
$ws.Range("A14").Value = 0.0035014989164352528
$ws.Range("B14").Value = -0.0035017689938721119
$ws.Range("A15").Value = -0.0080524472134415603
$ws.Range("B15").Value = 0.008034224780068655
$ws.Range("A16").Value = -0.0060342248612812455
$ws.Range("B16").Value = 0.0060032119353343738
$ws.Range("A17").Value = -0.0040032120180404362
$ws.Range("B17").Value = 0.003999999899774842
$ws.Range("A18").Value = -0.016101634074534132
$ws.Range("B18").Value = 0.016090817514808009
$ws.Range("A19").Value = -0.01209081755242547
$ws.Range("B19").Value = 0.012016033648514668
$ws.Range("A20").Value = -0.0080160336889409933
$ws.Range("B20").Value = 0.0080055446281441789
$ws.Range("A21").Value = -0.0040055446690150376
$ws.Range("B21").Value = 0.0039999999587747581
$ws.Range("A22").Value = -0.045707556231919355
$ws.Range("B22").Value = 0.045495858324899885
$ws.Range("A23").Value = -0.040495858391897954
$ws.Range("B23").Value = 0.040098282900659754
$ws.Range("A24").Value = -0.020098283101924075
$ws.Range("B24").Value = 0.019999999796017853
$ws.Range("A25").Value = -0.04161681314070087
$ws.Range("B25").Value = 0.041592170048858179
$ws.Range("A26").Value = -0.039092170119214842
$ws.Range("B26").Value = 0.039063216118211841
$ws.Range("A27").Value = -0.036563216189715142
$ws.Range("B27").Value = 0.036402539930092459
$ws.Range("A28").Value = -0.034402540004005999
$ws.Range("B28").Value = 0.034307123524540906
$ws.Range("A29").Value = -0.027307123645699427
$ws.Range("B29").Value = 0.027290002980227435
$ws.Range("A30").Value = 0.032709996443627265
$ws.Range("B30").Value = -0.032779378393593905
$ws.Range("A31").Value = 0.039779378275424548
$ws.Range("B31").Value = -0.039823170188501678
$ws.Range("A32").Value = -0.0040008795237191919
$ws.Range("B32").Value = 0.0039999999092437122
